$d = $word.ActiveDocument

# Grab the existing bullet-list template (numId 4 in the OOXML) that is
# already used elsewhere in the document, so the new bullet item below can
# continue that same list instead of starting a brand new numbering
# definition.
$listAnchor = $d.Content
$listAnchor.Find.Execute("n: #of iterations used to build and evaluate the classifier", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$existingListTemplate = $listAnchor.Paragraphs(1).Range.ListFormat.ListTemplate

# Locate the paragraph that ends with the "(sat.arff amd wine-white.arff)" sentence.
$anchorText = "The time needed to train and test the same data set is recorded for each algorithm and plotted for the two sets (sat.arff amd wine-white.arff)"

$r = $d.Content
$r.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# After Execute, $r is collapsed to the matched text. Grab the paragraph that
# contains it, then collapse a range to the very end of that paragraph (after
# its paragraph mark) so new paragraphs get inserted right after it.
$hostParagraph = $r.Paragraphs(1)
$hostRange = $hostParagraph.Range
$hostRange.Collapse(0)  # wdCollapseEnd

# Insert the first new paragraph: "Parameters to run:"
$hostRange.InsertParagraphAfter()
$hostRange.InsertAfter("Parameters to run:")
$p1Range = $hostRange.Paragraphs(1).Range
$p1Range.Style = "List Paragraph"
$p1Range.Font.Size = 12
$p1Range.Font.SizeBi = 12

# Move past the paragraph we just created.
$hostRange.Collapse(0)
$hostRange.Move(1, 1) | Out-Null

# Insert the second new paragraph: the "file:" bullet item.
$fileText = "file: data set file; if a value different than sat.arff or wine-white.arff is provided, the classifiers used have default parameters, instead of " + [char]8220 + "best" + [char]8221 + " ones."

$hostRange.InsertParagraphAfter()
$hostRange.InsertAfter($fileText)
$p2Range = $hostRange.Paragraphs(1).Range
$p2Range.Style = "List Paragraph"
$p2Range.Font.Size = 12
$p2Range.Font.SizeBi = 12
$p2Range.ListFormat.ApplyListTemplateWithLevel($existingListTemplate, $true, 1, $false)
